# Update "想去人数" (want-to-go count, column F) figures on each sheet
# to the refreshed numbers from the latest data pull.
$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 2273
$ws.Range("F3").Value = 336
$ws.Range("F4").Value = 178
$ws.Range("F5").Value = 182
$ws.Range("F6").Value = 342
$ws.Range("F8").Value = 704
$ws.Range("F9").Value = 514
$ws.Range("F10").Value = 665
$ws.Range("F11").Value = 370
$ws.Range("F13").Value = 359
$ws.Range("F14").Value = 972
$ws.Range("F15").Value = 3757
$ws.Range("F16").Value = 134
$ws.Range("F17").Value = 14
$ws.Range("F18").Value = 22
$ws.Range("F19").Value = 250
$ws.Range("F20").Value = 140
$ws.Range("F21").Value = 108
$ws.Range("F23").Value = 85
$ws.Range("F25").Value = 260
$ws.Range("F26").Value = 105

# Sheet "演出" (Performances)
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 66
$ws.Range("F3").Value = 33
$ws.Range("F6").Value = 178
$ws.Range("F7").Value = 209
$ws.Range("F8").Value = 2739
$ws.Range("F14").Value = 106
$ws.Range("F16").Value = 2490

# Sheet "本地生活" (Local Life)
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 257
$ws.Range("F3").Value = 40
$ws.Range("F4").Value = 394
$ws.Range("F5").Value = 172

# Sheet "全部类型" (All Types)
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 257
$ws.Range("F3").Value = 40
$ws.Range("F4").Value = 66
$ws.Range("F5").Value = 33
$ws.Range("F6").Value = 2273
$ws.Range("F7").Value = 394
$ws.Range("F8").Value = 336
$ws.Range("F9").Value = 178
$ws.Range("F10").Value = 182
$ws.Range("F11").Value = 342
$ws.Range("F15").Value = 178
$ws.Range("F16").Value = 172
$ws.Range("F17").Value = 704
$ws.Range("F18").Value = 514
$ws.Range("F19").Value = 665
$ws.Range("F20").Value = 370
$ws.Range("F22").Value = 359
$ws.Range("F23").Value = 972
$ws.Range("F24").Value = 3767
$ws.Range("F25").Value = 209
$ws.Range("F26").Value = 2739
$ws.Range("F30").Value = 134
$ws.Range("F31").Value = 14
$ws.Range("F32").Value = 22
$ws.Range("F35").Value = 250
$ws.Range("F36").Value = 140
$ws.Range("F37").Value = 108
$ws.Range("F39").Value = 106
$ws.Range("F41").Value = 85
$ws.Range("F43").Value = 260
$ws.Range("F44").Value = 105
$ws.Range("F45").Value = 2490
